$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# summer 24 week 15 inputs
$ws.Range("G5").Value = 0.79
$ws.Range("E6").Value = 1.33
$ws.Range("E7").Value = 1.86
$ws.Range("F7").Value = 1.48
